$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Turn off the existing (stale) AutoFilter before editing so the range/name
# get recomputed cleanly once we reapply it below.
$ws.AutoFilterMode = $false

# Row 83 is a "costaria" entry with an obvious length typo (length_cm = 11.5,
# far out of line with its mass_g = 121.5 and all other costaria rows).
# Delete it outright; every row below shifts up by one.
$ws.Rows.Item(83).Delete()

# Re-apply the AutoFilter across the full (now 121-row) data range, filtering
# species (column D) to "costaria" and visible_sori (column J) to "yes".
$rng = $ws.Range("A1:L121")
$rng.AutoFilter(4, "costaria", 7)
$rng.AutoFilter(10, "yes", 7)

# Keep the workbook-level _FilterDatabase defined name in sync with the new
# filter range.
$fd = $wb.Names.Item("Sheet1!_FilterDatabase")
$fd.RefersTo = "=Sheet1!`$A`$1:`$L`$121"

# Select the first visible (filtered-in) data row, matching the author's
# final selection after filtering.
$ws.Rows.Item(83).Select()
